$wb = $excel.ActiveWorkbook

# Add the new "TImers" sheet (initially inserted before the active sheet)
$ws = $wb.Worksheets.Add()
$ws.Name = "TImers"

# Populate cell labels in the same order the original author typed them in,
# so the shared-strings table comes out in the matching order.
$ws.Range("A1").Value = "Desired interrupt period"
$ws.Range("A6").Value = "Tcnt"
$ws.Range("A5").Value = "Tarr"
$ws.Range("A3").Value = "SysClk"
$ws.Range("A4").Value = "Tpsc"
$ws.Range("A2").Value = "Actual Period"

# Values / formulas
$ws.Range("B1").Value = 0.0001
$ws.Range("B2").Formula = "=B6*B5"
$ws.Range("B3").Value = 84000000
$ws.Range("B4").Value = 0
$ws.Range("C4").Formula = "=((B1*B3)/B5)-1"
$ws.Range("B5").Value = 8400
$ws.Range("C5").Formula = "=B1/B6"
$ws.Range("E5").Formula = "=4000/0.0001"
$ws.Range("B6").Formula = "=(B4+1)/B3"
$ws.Range("E6").Formula = "=B6*12"

# Number formats for the high-precision timing cells
$ws.Range("B2").NumberFormat = "0.00000000"
$ws.Range("B6").NumberFormat = "0.0000000000"
$ws.Range("E6").NumberFormat = "0.000000000000"

# Move the new sheet to the end (after CurrentSensingOpAmpGain)
$ws.Move($null, $wb.Worksheets.Item("CurrentSensingOpAmpGain"))

# Re-fetch a live reference post-move, select E13, and make it the active tab
$tImers = $wb.Worksheets.Item("TImers")
$tImers.Range("E13").Select()
$tImers.Activate()
